# Daily "refresh" update for the tracking sheet.
# For every data row:
#   - D = total number of days in the current cycle
#   - E = number of days remaining in the current cycle
#   - F = start date (yyyyMMdd) of the current cycle
#
# The "today" reference date advances by one day (to 2025-11-03). For each
# row we recompute how many days remain until the cycle (start date + total
# days) ends, as of the new "today":
#   remaining = (F_old + D) - today
# If the cycle has ended (remaining <= 0) a new cycle starts today:
#   E_new = D ; F_new = today
# Otherwise the cycle continues with the same start date, just fewer days
# remaining:
#   E_new = remaining ; F_new = F_old

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = [DateTime]::ParseExact("20251103", "yyyyMMdd", $null)
$todayOA = $today.ToOADate()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value()
    $fVal = $fCell.Value()

    if ($dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    try {
        $fStr = [string]([int64]$fVal)
        $fDate = [DateTime]::ParseExact($fStr, "yyyyMMdd", $null)
    } catch {
        # Malformed start date (e.g. "202510929") - leave the row untouched,
        # matching the reference workbook's behaviour.
        continue
    }

    $dInt = [int64]$dVal
    $endOA = $fDate.ToOADate() + $dInt
    $remaining = $endOA - $todayOA

    if ($remaining -le 0) {
        $eCell.Value = $dInt
        $fCell.Value = [int64]($today.ToString("yyyyMMdd"))
    } else {
        $eCell.Value = [int64]$remaining
    }
}
